# Updated cryptos list on Thu Sep 28 21:10:14 UTC 2023 with GitHub Actions
#
# This script reproduces the upstream crypto-price refresh: for each coin
# row, the "Price" (column D) and "Volume(1h)" (column E) cells are updated
# to the latest scraped figures. The cells hold plain text (not numbers),
# so we temporarily force a text NumberFormat before writing the value -
# this stops Excel from re-interpreting values such as "215.10" or
# "27.070.23" as numbers (which would silently drop meaningful trailing
# zeros / thousands separators) - and then restore the cell to the
# workbook's default "Normal" style so no visible formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Address, $Text) {
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value2 = $Text
    $cell.Style = "Normal"
}

Set-TextCell "D2" "27.070.23"
Set-TextCell "E2" "  +3.10%  "
Set-TextCell "D3" "1.656.92"
Set-TextCell "E3" "  +3.90%  "
Set-TextCell "E4" "  +0.08%  "
Set-TextCell "D5" "215.10"
Set-TextCell "E5" "  +1.84%  "
Set-TextCell "E6" "  +0.86%  "
Set-TextCell "E7" "  +0.02%  "
Set-TextCell "E8" "  +1.80%  "
Set-TextCell "E9" "  +1.72%  "
Set-TextCell "D10" "19.70"
Set-TextCell "E10" "  +3.91%  "
Set-TextCell "D11" "0.0864"
Set-TextCell "E11" "  +1.04%  "
Set-TextCell "D12" "1.890.81"
Set-TextCell "E12" "  +3.94%  "
Set-TextCell "D13" "1.661.62"
Set-TextCell "E13" "  +4.38%  "
Set-TextCell "D14" "4.07"
Set-TextCell "E14" "  +2.20%  "
Set-TextCell "E15" "  +3.04%  "
Set-TextCell "D16" "65.20"
Set-TextCell "E16" "  +2.67%  "
Set-TextCell "D17" "27.080.80"
Set-TextCell "E17" "  +3.15%  "
Set-TextCell "D18" "238.99"
Set-TextCell "E18" "  +3.77%  "
Set-TextCell "D19" "7.85"
Set-TextCell "E19" "  +2.80%  "
Set-TextCell "E20" "  +1.35%  "
Set-TextCell "E21" "  -0.01%  "
Set-TextCell "E22" "  +4.52%  "
Set-TextCell "E23" "  +4.30%  "
Set-TextCell "E24" "  +3.64%  "
Set-TextCell "D25" "145.94"
Set-TextCell "E25" "  -0.08%  "
Set-TextCell "E26" "  -0.07%  "
Set-TextCell "D27" "7.13"
Set-TextCell "E27" "  +2.21%  "
Set-TextCell "E28" "  +1.44%  "
Set-TextCell "D29" "15.85"
Set-TextCell "E29" "  +3.35%  "
Set-TextCell "E30" "  +0.79%  "
Set-TextCell "E31" "  +1.90%  "
Set-TextCell "D32" "3.30"
Set-TextCell "E32" "  +3.24%  "
Set-TextCell "D33" "1.518.04"
Set-TextCell "E33" "  +3.36%  "
Set-TextCell "D34" "3.07"
Set-TextCell "E34" "  +4.55%  "
Set-TextCell "D35" "1.60"
Set-TextCell "E35" "  +9.73%  "
Set-TextCell "E36" "  +0.03%  "
Set-TextCell "D37" "0.576"
Set-TextCell "E37" "  +1.81%  "
Set-TextCell "E38" "  +8.60%  "
Set-TextCell "E40" "  +3.78%  "
Set-TextCell "E41" "  +0.01%  "
Set-TextCell "E42" "  +4.34%  "
Set-TextCell "D43" "66.02"
Set-TextCell "E43" "  +9.40%  "
Set-TextCell "D44" "1.798.27"
Set-TextCell "E44" "  +3.77%  "
Set-TextCell "D45" "0.777"
Set-TextCell "E45" "  +3.16%  "
Set-TextCell "E46" "  -1.13%  "
Set-TextCell "D47" "89.57"
Set-TextCell "E47" "  +1.80%  "
Set-TextCell "D48" "0.0₆0106"
Set-TextCell "E48" "  +0.69%  "
Set-TextCell "E49" "  +3.46%  "
Set-TextCell "E50" "  +0.88%  "
Set-TextCell "E51" "  +2.96%  "
